# Add two new sheets ("linear" and "polynomial") after the existing
# "Sheet1", populate them with ultrasonic-sensor reading data, resize
# their columns, and move the active-sheet/selection state around to
# match the edited workbook.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Create the two new worksheets, right after Sheet1, in order.
# ---------------------------------------------------------------------
$linear = $wb.Worksheets.Add([System.Type]::Missing, $sheet1)
$linear.Name = "linear"

$poly = $wb.Worksheets.Add([System.Type]::Missing, $linear)
$poly.Name = "polynomial"

# ---------------------------------------------------------------------
# "linear" sheet data
# ---------------------------------------------------------------------
$linear.Range("A1").Value = "Speed:"
$linear.Range("B1").Value = 340
$linear.Range("A2").Value = "Measurement on Tape (cm)"
$linear.Range("B2").Value = "Actual Readings from sensor) cm"

$linear.Range("A3").Value = 10
$linear.Range("B3").Value = 9.1300000000000008
$linear.Range("A4").Value = 14
$linear.Range("B4").Value = 13.5
$linear.Range("A5").Value = 22
$linear.Range("B5").Value = 20
$linear.Range("A6").Value = 27
$linear.Range("B6").Value = 24.8
$linear.Range("A7").Value = 31
$linear.Range("B7").Value = 29.5
$linear.Range("A8").Value = 35
$linear.Range("B8").Value = 32.299999999999997
$linear.Range("A9").Value = 40
$linear.Range("B9").Value = 37.6

$linear.Columns.Item(1).ColumnWidth = 36.944010416666664
$linear.Columns.Item(2).ColumnWidth = 33.166666666666664

# ---------------------------------------------------------------------
# "polynomial" sheet data
# ---------------------------------------------------------------------
$poly.Range("A1").Value = "Speed:"
$poly.Range("B1").Value = 340
$poly.Range("A2").Value = "Measurement on Tape (cm)"
$poly.Range("B2").Value = "Actual Readings from sensor) cm"

$poly.Range("A3").Value = 10
$poly.Range("B3").Value = 10.1
$poly.Range("A4").Value = 14
$poly.Range("B4").Value = 14.46
$poly.Range("A5").Value = 22
$poly.Range("B5").Value = 20.5
$poly.Range("A6").Value = 27
$poly.Range("B6").Value = 25.36
$poly.Range("A7").Value = 31
$poly.Range("B7").Value = 28.7
$poly.Range("A8").Value = 35
$poly.Range("B8").Value = 32.270000000000003
$poly.Range("A9").Value = 40
$poly.Range("B9").Value = 37

$poly.Columns.Item(1).ColumnWidth = 36.944010416666664
$poly.Columns.Item(2).ColumnWidth = 33.166666666666664

# ---------------------------------------------------------------------
# Selections on each sheet (Excel records the active cell/selection
# per sheet, plus which sheet tab is active overall). "polynomial" is
# selected last so it ends up the active tab.
# ---------------------------------------------------------------------
[void]$sheet1.Select()
[void]$sheet1.Range("A2:B3").Select()

[void]$linear.Select()
[void]$linear.Range("A26").Select()

[void]$poly.Select()
[void]$poly.Range("G18").Select()
